# "Beschrijving aangepast naar 3.0"
# Reworks the "Berekening voor nieuwe klas" section of Blad1: adds two new
# input rows (grams per pack of hagelslag / per tub of margarine) and makes
# the "doosjes"/"kuipjes" labels computed (CONCAT) instead of static text,
# and bumps the class size (leerlingen + begeleiders) from 30 to 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Make room: insert 3 blank rows right before the old "Berekening voor
# nieuwe klas" block (old row 34). Excel shifts every formula below along
# with the rows, so A35->A38, A36->A39, ..., A41->A44 automatically, and
# all formulas referencing them are auto-rewritten too.
$ws.Rows.Item(30).Resize(3).Insert()

# New block: "Hagelslag" (grams per pack)
$ws.Cells.Item(30, 1).Value = "Hagelslag"
$ws.Cells.Item(30, 1).Font.Bold = $true
$ws.Cells.Item(30, 1).Font.Size = 14
$ws.Rows.Item(30).RowHeight = 18.75

$ws.Cells.Item(31, 1).Value = 390
$ws.Cells.Item(31, 2).Value = "g/pak"

# New block: "Margarine" (grams per tub)
$ws.Cells.Item(32, 1).Value = "Margarine"
$ws.Cells.Item(32, 1).Font.Bold = $true
$ws.Cells.Item(32, 1).Font.Size = 14
$ws.Rows.Item(32).RowHeight = 18.75

$ws.Cells.Item(33, 1).Value = 250
$ws.Cells.Item(33, 2).Value = "g/kuipje"

# "Berekening voor nieuwe klas" block (now rows 37-44): update class size
$ws.Cells.Item(38, 1).Value = 47

# Row 40 ("kuipjes van ... gram"): divide by the new margarine-tub-size
# input cell (A33) instead of the hardcoded 250, and compute the label
# text from that cell too.
$ws.Cells.Item(40, 1).Formula = "=ROUNDUP(A39/A33,0)"
$ws.Cells.Item(40, 2).Formula = '=CONCAT("kuipjes van ",A33," gram")'

# Row 42 ("doosjes hagelslag van ... gram"): divide by the new
# hagelslag-pack-size input cell (A31) instead of the hardcoded 250, and
# compute the label text from that cell too.
$ws.Cells.Item(42, 1).Formula = "=ROUNDUP(A41/A31,0)"
$ws.Cells.Item(42, 2).Formula = '=CONCAT("doosjes hagelslag van ",A31," gram")'

# Match the author's final selection/scroll position.
[void]$ws.Range("H36").Select()
